# Update odds values in Sheet1 per the 2026-01-07 data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.46
$ws.Range("G2").Value = 2.72
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.85
$ws.Range("L2").Value = 1.6
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 2.54
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 1.52
$ws.Range("Q2").Value = 2.58
$ws.Range("T2").Value = 2.04
$ws.Range("U2").Value = 1.79
$ws.Range("V2").Value = 1.36
$ws.Range("W2").Value = 1.58
$ws.Range("X2").Value = 9
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 29
$ws.Range("AB2").Value = 9
$ws.Range("AD2").Value = 16.5
$ws.Range("AF2").Value = 18
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 25
$ws.Range("AK2").Value = 1000

# Row 3
$ws.Range("F3").Value = 3.25
$ws.Range("G3").Value = 3.7
$ws.Range("H3").Value = 2.06
$ws.Range("L3").Value = 1.32
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 4.7
$ws.Range("O3").Value = 1.21
$ws.Range("Q3").Value = 1.63
$ws.Range("S3").Value = 2.58
$ws.Range("Y3").Value = 15.5
$ws.Range("AA3").Value = 30
$ws.Range("AB3").Value = 20
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 22
$ws.Range("AF3").Value = 32
$ws.Range("AG3").Value = 18.5
$ws.Range("AH3").Value = 18
$ws.Range("AI3").Value = 30
$ws.Range("AJ3").Value = 60
$ws.Range("AK3").Value = 38
$ws.Range("AL3").Value = 46
$ws.Range("AM3").Value = 75
$ws.Range("AN3").Value = 30
$ws.Range("AO3").Value = 12.5

# Row 4
$ws.Range("F4").Value = 4.8
$ws.Range("G4").Value = 5.8
$ws.Range("H4").Value = 1.81
$ws.Range("I4").Value = 1.97
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.85
$ws.Range("L4").Value = 1.46
$ws.Range("N4").Value = 3.05
$ws.Range("O4").Value = 1.4
$ws.Range("S4").Value = 4
$ws.Range("Y4").Value = 8.8
$ws.Range("AA4").Value = 25
$ws.Range("AB4").Value = 18
$ws.Range("AC4").Value = 9.6
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 40
$ws.Range("AG4").Value = 25
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 1000
$ws.Range("AK4").Value = 95
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 19.5

# Row 5
$ws.Range("F5").Value = 2.66
$ws.Range("G5").Value = 2.98
$ws.Range("H5").Value = 2.28
$ws.Range("I5").Value = 2.68
$ws.Range("J5").Value = 3.65
$ws.Range("K5").Value = 4.4
$ws.Range("L5").Value = 1.26
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 3.45
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 2.24
$ws.Range("Q5").Value = 1.53
$ws.Range("R5").Value = 1.55
$ws.Range("S5").Value = 2.26
$ws.Range("T5").Value = 1.5
$ws.Range("U5").Value = 2.44
$ws.Range("V5").Value = 1.6
$ws.Range("W5").Value = 1.51
$ws.Range("X5").Value = 26
$ws.Range("Y5").Value = 17
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 55
$ws.Range("AB5").Value = 22
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 13.5
$ws.Range("AE5").Value = 25
$ws.Range("AF5").Value = 25
$ws.Range("AG5").Value = 14
$ws.Range("AH5").Value = 15.5
$ws.Range("AI5").Value = 32
$ws.Range("AJ5").Value = 55
$ws.Range("AK5").Value = 28
$ws.Range("AL5").Value = 32
$ws.Range("AM5").Value = 60
$ws.Range("AN5").Value = 19.5
$ws.Range("AO5").Value = 14

# Row 6
$ws.Range("F6").Value = 1.96
$ws.Range("G6").Value = 2.68
$ws.Range("H6").Value = 1.09
$ws.Range("I6").Value = 4.3
$ws.Range("K6").Value = 5.7
$ws.Range("L6").Value = 1.3
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 3.05
$ws.Range("O6").Value = 1.21
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.56
$ws.Range("R6").Value = 1.46
$ws.Range("S6").Value = 2.4
$ws.Range("T6").Value = 1.56
$ws.Range("U6").Value = 2.28
$ws.Range("V6").Value = 1.33
$ws.Range("W6").Value = 1.6
$ws.Range("Y6").Value = 24
$ws.Range("Z6").Value = 38
$ws.Range("AB6").Value = 18.5
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 21
$ws.Range("AE6").Value = 50
$ws.Range("AF6").Value = 24
$ws.Range("AG6").Value = 17
$ws.Range("AH6").Value = 23
$ws.Range("AJ6").Value = 44
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 46
$ws.Range("AN6").Value = 19
$ws.Range("AO6").Value = 38

# Row 7
$ws.Range("F7").Value = 1.66
$ws.Range("G7").Value = 2.18
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 5.4
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 1.29
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 2.96
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.36
$ws.Range("R7").Value = 1.35
$ws.Range("S7").Value = 2.02
$ws.Range("T7").Value = 1.44
$ws.Range("U7").Value = 1.89
$ws.Range("V7").Value = 1.23
$ws.Range("W7").Value = 1.07
$ws.Range("X7").Value = 30
$ws.Range("Y7").Value = 29
$ws.Range("Z7").Value = 50
$ws.Range("AB7").Value = 17
$ws.Range("AC7").Value = 13.5
$ws.Range("AD7").Value = 25
$ws.Range("AF7").Value = 20
$ws.Range("AG7").Value = 15.5
$ws.Range("AH7").Value = 25
$ws.Range("AJ7").Value = 30
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 44
$ws.Range("AN7").Value = 14.5

# Row 8
$ws.Range("J8").Value = 5.2
$ws.Range("K8").Value = 5.3
$ws.Range("S8").Value = 3.8
$ws.Range("T8").Value = 2.54
$ws.Range("U8").Value = 1.61
$ws.Range("X8").Value = 15.5
$ws.Range("Y8").Value = 29
$ws.Range("AC8").Value = 12
$ws.Range("AD8").Value = 44
$ws.Range("AF8").Value = 6.4
$ws.Range("AK8").Value = 16.5

# Row 9
$ws.Range("G9").Value = 3.1
$ws.Range("U9").Value = 2.2

# Row 10
$ws.Range("F10").Value = 3.6
$ws.Range("G10").Value = 3.65
$ws.Range("J10").Value = 3.8
$ws.Range("K10").Value = 3.85
$ws.Range("W10").Value = 1.37
$ws.Range("AF10").Value = 26
$ws.Range("AL10").Value = 46

# Row 11
$ws.Range("I11").Value = 4.9
$ws.Range("AB11").Value = 8.2
$ws.Range("AI11").Value = 75
$ws.Range("AM11").Value = 100

# Row 12
$ws.Range("F12").Value = 3.3
$ws.Range("G12").Value = 3.35
$ws.Range("H12").Value = 2.38
$ws.Range("I12").Value = 2.4
$ws.Range("P12").Value = 2
$ws.Range("V12").Value = 1.71
$ws.Range("W12").Value = 1.42
$ws.Range("AB12").Value = 13
$ws.Range("AF12").Value = 23
$ws.Range("AG12").Value = 14
$ws.Range("AO12").Value = 19.5

# Row 13
$ws.Range("X13").Value = 17
$ws.Range("Y13").Value = 14.5
$ws.Range("AD13").Value = 14.5
$ws.Range("AO13").Value = 32

# Row 14
$ws.Range("O14").Value = 1.15
$ws.Range("Q14").Value = 1.44
$ws.Range("R14").Value = 1.88
$ws.Range("S14").Value = 2.1
$ws.Range("T14").Value = 1.65
$ws.Range("U14").Value = 2.48
$ws.Range("AH14").Value = 18.5
$ws.Range("AJ14").Value = 13.5
$ws.Range("AO14").Value = 75

# Row 15
$ws.Range("Z15").Value = 36
$ws.Range("AD15").Value = 19.5
$ws.Range("AF15").Value = 10.5
$ws.Range("AM15").Value = 140

# Row 16
$ws.Range("L16").Value = 1.5
$ws.Range("X16").Value = 9.6

# Row 17
$ws.Range("K17").Value = 3.15
$ws.Range("P17").Value = 1.68
$ws.Range("Q17").Value = 2.42
$ws.Range("AH17").Value = 20

# Row 18
$ws.Range("H18").Value = 1.41
$ws.Range("I18").Value = 1.42
$ws.Range("J18").Value = 5.2
$ws.Range("K18").Value = 5.3
$ws.Range("N18").Value = 4.7
$ws.Range("Q18").Value = 1.77
$ws.Range("T18").Value = 2.04
$ws.Range("V18").Value = 3.4
$ws.Range("W18").Value = 1.11
$ws.Range("Z18").Value = 8
$ws.Range("AJ18").Value = 360
$ws.Range("AL18").Value = 130

# Row 19
$ws.Range("G19").Value = 5.4
$ws.Range("N19").Value = 4.4
$ws.Range("Q19").Value = 1.81
$ws.Range("S19").Value = 3.05
$ws.Range("Y19").Value = 9.4
$ws.Range("AH19").Value = 18.5
$ws.Range("AN19").Value = 80
$ws.Range("AO19").Value = 9

# Row 20
$ws.Range("G20").Value = 1.75
$ws.Range("N20").Value = 4.3
$ws.Range("T20").Value = 1.84
$ws.Range("W20").Value = 2.32
$ws.Range("X20").Value = 16.5
$ws.Range("AB20").Value = 9.2
$ws.Range("AE20").Value = 65
$ws.Range("AJ20").Value = 17
$ws.Range("AN20").Value = 9.8
